# Auto-generated update of Golem_Profits leve-profit figures (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 586.25
$ws.Range("I33").Value = 589.0833
$ws.Range("K33").Value = 589.0833
$ws.Range("M33").Value = -360.0833

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2016.3
$ws.Range("I38").Value = 18.11111
$ws.Range("J38").Value = 20000
$ws.Range("K38").Value = 54.33333
$ws.Range("L38").Value = 60000
$ws.Range("M38").Value = 317.66667
$ws.Range("N38").Value = -60744

# ALC row 42
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 585.1429000000001
$ws.Range("I42").Value = 233.33333
$ws.Range("K42").Value = 699.99999
$ws.Range("M42").Value = -469.99999

# ALC row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 250.6
$ws.Range("J53").Value = 359.33334
$ws.Range("L53").Value = 359.33334
$ws.Range("N53").Value = -1633.33334

# ALC row 75
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 49997.5
$ws.Range("J75").Value = 49997.5
$ws.Range("L75").Value = 49997.5
$ws.Range("N75").Value = -51869.5

# ALC row 78
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 49997.5
$ws.Range("J78").Value = 49997.5
$ws.Range("L78").Value = 149992.5
$ws.Range("N78").Value = -159352.5

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2520.4211
$ws.Range("I138").Value = 713.6667
$ws.Range("K138").Value = 2141.0001
$ws.Range("M138").Value = 2998.9999

# ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 612
$ws.Range("I4").Value = 612
$ws.Range("K4").Value = 612
$ws.Range("M4").Value = -496

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2163.1
$ws.Range("I45").Value = 2163.1
$ws.Range("K45").Value = 2163.1
$ws.Range("M45").Value = -1786.1

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 179
$ws.Range("I7").Value = 166.33333
$ws.Range("K7").Value = 166.33333
$ws.Range("M7").Value = -53.33332999999999

# CRP row 38
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 26200
$ws.Range("I38").Value = 24500
$ws.Range("J38").Value = 27333.334
$ws.Range("K38").Value = 24500
$ws.Range("L38").Value = 27333.334
$ws.Range("M38").Value = -24123
$ws.Range("N38").Value = -28087.334

# CRP row 44
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35884

# CRP row 46
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 26200
$ws.Range("I46").Value = 24500
$ws.Range("J46").Value = 27333.334
$ws.Range("K46").Value = 24500
$ws.Range("L46").Value = 27333.334
$ws.Range("M46").Value = -24289
$ws.Range("N46").Value = -27755.334

# CRP row 125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -34920

# CUL row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 711.2
$ws.Range("I6").Value = 52
$ws.Range("J6").Value = 1700
$ws.Range("K6").Value = 156
$ws.Range("L6").Value = 5100
$ws.Range("M6").Value = -43
$ws.Range("N6").Value = -5326

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 321.6
$ws.Range("J17").Value = 1500
$ws.Range("L17").Value = 4500
$ws.Range("N17").Value = -4838

# CUL row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 3860.2856
$ws.Range("I99").Value = 4403.6665
$ws.Range("J99").Value = 600
$ws.Range("K99").Value = 13210.9995
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = -10964.9995
$ws.Range("N99").Value = -6292

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 147.24
$ws.Range("I2").Value = 143.28572
$ws.Range("J2").Value = 152.27272
$ws.Range("K2").Value = 143.28572
$ws.Range("L2").Value = 152.27272
$ws.Range("M2").Value = -30.28572
$ws.Range("N2").Value = -378.27272

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 63000
$ws.Range("I5").Value = 63000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 63000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -62888
$ws.Range("N5").Value = ""

# GSM row 9
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 221.35715
$ws.Range("I9").Value = 145.36363
$ws.Range("K9").Value = 145.36363
$ws.Range("M9").Value = 24.63637

# GSM row 17
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 1500
$ws.Range("J17").Value = 1500
$ws.Range("L17").Value = 1500
$ws.Range("N17").Value = -1836

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 800
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4006.4375
$ws.Range("I68").Value = 4355.222
$ws.Range("J68").Value = 3558
$ws.Range("K68").Value = 4355.222
$ws.Range("L68").Value = 3558
$ws.Range("M68").Value = -3606.222
$ws.Range("N68").Value = -5056

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4006.4375
$ws.Range("I71").Value = 4355.222
$ws.Range("J71").Value = 3558
$ws.Range("K71").Value = 21776.11
$ws.Range("L71").Value = 17790
$ws.Range("M71").Value = -18032.11
$ws.Range("N71").Value = -25278

# WVR row 21
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 200
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 200
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = -670

# WVR row 30
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = ""

# WVR row 35
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 200
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 200
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 200
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = -780

# WVR row 55
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 26236.363
$ws.Range("I55").Value = 8833.333000000001
$ws.Range("J55").Value = 32762.5
$ws.Range("K55").Value = 8833.333000000001
$ws.Range("L55").Value = 32762.5
$ws.Range("M55").Value = -8556.333000000001
$ws.Range("N55").Value = -33316.5

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2025
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -11650

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2146
$ws.Range("I136").Value = 2375
$ws.Range("K136").Value = 7125
$ws.Range("M136").Value = -4575
